$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F16 value from 20 to 24 (dependent formulas in G16/H16 recalc automatically)
$ws.Range("F16").Value = 24

# Move the active cell selection from F17 to F16
$ws.Range("F16").Select()
